{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer cells in the (only)\n// table with their new values, in document order. The mapping is\n// strictly positional (the Nth populated cell gets the Nth new\n// value) -- it is NOT a text find/replace, since some new values are\n// identical to *other* cells' old values (e.g. cell(0,0) old text\n// \"17\u00f76=2, 5\" also becomes the new text of cell(4,0)).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// New answers, in the same left-to-right / top-to-bottom order as the\n// cells that currently hold an answer (every 4th row starting at 0).\nconst newAnswers = [\n  \"46\u00f73=15, 1\", \"68\u00f77=9, 5\", \"86\u00f78=10, 6\", \"43\u00f78=5, 3\", \"19\u00f75=3, 4\",\n  \"17\u00f76=2, 5\", \"21\u00f76=3, 3\", \"67\u00f79=7, 4\", \"88\u00f72=44, 0\", \"98\u00f74=24, 2\",\n  \"26\u00f74=6, 2\", \"80\u00f78=10, 0\", \"10\u00f73=3, 1\", \"30\u00f74=7, 2\", \"23\u00f78=2, 7\",\n  \"39\u00f77=5, 4\", \"13\u00f79=1, 4\", \"14\u00f75=2, 4\", \"81\u00f73=27, 0\", \"46\u00f74=11, 2\",\n  \"14\u00f78=1, 6\", \"47\u00f79=5, 2\", \"61\u00f72=30, 1\", \"63\u00f75=12, 3\", \"21\u00f77=3, 0\",\n];\n\nconst rowValues = table.values;\nlet answerIndex = 0;\nfor (let r = 0; r < rowValues.length; r++) {\n  const row = rowValues[r];\n  // A data row here has non-empty text in its first cell; blank\n  // spacer rows have empty strings in every cell.\n  if (!row[0]) {\n    continue;\n  }\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newAnswers[answerIndex], Word.InsertLocation.replace);\n    answerIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer cells in the (only)\n# table with their new values, in document order. The mapping is\n# strictly positional (the Nth populated cell gets the Nth new\n# value) -- it is NOT a text find/replace, since some new values are\n# identical to *other* cells' old values (e.g. row1/col1's old text\n# \"17\u00f76=2, 5\" also becomes the new text of row5/col1).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# New answers, in the same left-to-right / top-to-bottom order as the\n# cells that currently hold an answer (every 4th row starting at\n# row 1, 1-indexed).\n$newAnswers = @(\n  \"46\u00f73=15, 1\", \"68\u00f77=9, 5\", \"86\u00f78=10, 6\", \"43\u00f78=5, 3\", \"19\u00f75=3, 4\",\n  \"17\u00f76=2, 5\", \"21\u00f76=3, 3\", \"67\u00f79=7, 4\", \"88\u00f72=44, 0\", \"98\u00f74=24, 2\",\n  \"26\u00f74=6, 2\", \"80\u00f78=10, 0\", \"10\u00f73=3, 1\", \"30\u00f74=7, 2\", \"23\u00f78=2, 7\",\n  \"39\u00f77=5, 4\", \"13\u00f79=1, 4\", \"14\u00f75=2, 4\", \"81\u00f73=27, 0\", \"46\u00f74=11, 2\",\n  \"14\u00f78=1, 6\", \"47\u00f79=5, 2\", \"61\u00f72=30, 1\", \"63\u00f75=12, 3\", \"21\u00f77=3, 0\"\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$answerIndex = 0\nforeach ($r in $dataRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newAnswers[$answerIndex]\n        $answerIndex++\n    }\n}\n"}
